$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet Hoja1, cell A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$text = $cellA1.Value2
$text = $text -replace [regex]::Escape("✅ 1000 Bs = 3.25 = 12494.84 pesos`n✅ 12494.84 pesos = 3.24 = 980.71 Bs"), "✅ 1000 Bs = 3.24 = 12349.51 pesos`n✅ 12349.51 pesos = 3.2 = 967.46 Bs"
$cellA1.Value2 = $text

# --- Update rate figures on sheet tasas ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 309
$wsTasas.Range("O10").Value = 3816
$wsTasas.Range("N12").Value = 3855
$wsTasas.Range("O12").Value = 302
